# Update country data (COVID-19 dataset) in sheet "Pais"
# Values derived from the target diff: rows whose data/labels changed position or content
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 17:05"

# Row 4
$ws.Range("B4").Value = 1411339
$ws.Range("C4").Value = 2703
$ws.Range("E4").Value = 1029132

# Row 9
$ws.Range("A9").Value = "Brasil"
$ws.Range("B9").Value = 179457
$ws.Range("C9").Value = 1855
$ws.Range("D9").Value = 72597
$ws.Range("E9").Value = 94329
$ws.Range("F9").Value = 8318
$ws.Range("G9").Value = 127
$ws.Range("H9").Value = 12531

# Row 10
$ws.Range("A10").Value = "Francia"
$ws.Range("B10").Value = 178225
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 57785
$ws.Range("E10").Value = 93449
$ws.Range("F10").Value = 2542
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 26991

# Row 30
$ws.Range("D30").Value = 4809
$ws.Range("E30").Value = 20516
$ws.Range("F30").Value = 19

# Row 61
$ws.Range("A61").Value = "Moldavia"
$ws.Range("B61").Value = 5406
$ws.Range("C61").Value = 252
$ws.Range("D61").Value = 2176
$ws.Range("E61").Value = 3045
$ws.Range("F61").Value = 251
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 185

# Row 62
$ws.Range("A62").Value = "Afganistan"
$ws.Range("B62").Value = 5226
$ws.Range("C62").Value = 263
$ws.Range("D62").Value = 648
$ws.Range("E62").Value = 4446
$ws.Range("F62").Value = 7
$ws.Range("G62").Value = 5
$ws.Range("H62").Value = 132

# Row 69
$ws.Range("A69").Value = "Irak"
$ws.Range("B69").Value = 3032
$ws.Range("C69").Value = 119
$ws.Range("D69").Value = 1966
$ws.Range("E69").Value = 951
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 115

# Row 70
$ws.Range("A70").Value = "Tailandia"
$ws.Range("B70").Value = 3017
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 2844
$ws.Range("E70").Value = 117
$ws.Range("F70").Value = 61
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 56

# Row 71
$ws.Range("A71").Value = "Bolivia"
$ws.Range("B71").Value = 2964
$ws.Range("C71").Value = 133
$ws.Range("D71").Value = 313
$ws.Range("E71").Value = 2523
$ws.Range("F71").Value = 3
$ws.Range("G71").Value = 6
$ws.Range("H71").Value = 128

# Row 104
$ws.Range("B104").Value = 893
$ws.Range("C104").Value = 4
$ws.Range("E104").Value = 502

# Row 142
$ws.Range("A142").Value = "Nepal"
$ws.Range("B142").Value = 243
$ws.Range("C142").Value = 26
$ws.Range("D142").Value = 35
$ws.Range("E142").Value = 208
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 0

# Row 143
$ws.Range("A143").Value = "Santo Tome y Principe"
$ws.Range("B143").Value = 220
$ws.Range("C143").Value = 12
$ws.Range("D143").Value = 4
$ws.Range("E143").Value = 210
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 6

# Row 146
$ws.Range("A146").Value = "Madagascar"
$ws.Range("B146").Value = 212
$ws.Range("C146").Value = 26
$ws.Range("D146").Value = 107
$ws.Range("E146").Value = 105
$ws.Range("F146").Value = 1
$ws.Range("H146").Value = 0

# Row 147
$ws.Range("A147").Value = "Togo"
$ws.Range("B147").Value = 199
$ws.Range("D147").Value = 92
$ws.Range("E147").Value = 96
$ws.Range("H147").Value = 11

# Row 148
$ws.Range("A148").Value = "Sudan del Sur"
$ws.Range("B148").Value = 194
$ws.Range("D148").Value = 2
$ws.Range("E148").Value = 192
$ws.Range("F148").Value = 0
$ws.Range("H148").Value = 0

# Row 149
$ws.Range("A149").Value = "Martinica"
$ws.Range("D149").Value = 91
$ws.Range("E149").Value = 82
$ws.Range("F149").Value = 4
$ws.Range("H149").Value = 14

# Row 150
$ws.Range("A150").Value = "Islas Feroe"
$ws.Range("B150").Value = 187
$ws.Range("D150").Value = 187
$ws.Range("E150").Value = 0
$ws.Range("F150").Value = 0

# Row 175
$ws.Range("B175").Value = 48
$ws.Range("C175").Value = 1
$ws.Range("E175").Value = 16

# Row 193
$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2

# Row 194
$ws.Range("A194").Value = "Nueva Caledonia"
$ws.Range("D194").Value = 18
$ws.Range("H194").Value = 0
